# Auto-generated data refresh for 悖论模拟干员名单作者版.xlsx (Sheet1)
# Updates maa:// pass-rate annotations and the "last updated" timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("L3").Value = "*maa://22880 (66.31), maa://20276 (85.45), *maa://22749 (72.73)"
$ws.Range("X3").Value = "maa://27396 (84.35), maa://27484 (96.26), maa://27480 (82.86)"
$ws.Range("T4").Value = "maa://32509 (97.17), maa://27295 (84.62), maa://22754 (90.41), *maa://21746 (55.81), *maa://31008 (78.57)"
$ws.Range("L5").Value = "*maa://22757 (77.14)"
$ws.Range("D6").Value = "maa://42407 (95.0)"
$ws.Range("T7").Value = "maa://21291 (84.78)"
$ws.Range("A8").Value = "更新日期：2025.01.18 13:16:09"
$ws.Range("AB9").Value = "maa://28711 (87.38), ***maa://22740 (5.77), **maa://39938 (46.15), **maa://27377 (42.86), ***maa://25174 (19.05), maa://40166 (95.0)"
$ws.Range("AB16").Value = "maa://26228 (95.6)"
$ws.Range("D17").Value = "maa://21624 (83.78)"
$ws.Range("AB18").Value = "maa://24393 (97.67)"
$ws.Range("AB19").Value = "*maa://30709 (63.94), *maa://36668 (56.41)"
$ws.Range("L20").Value = "maa://41331 (84.96)"
$ws.Range("AF21").Value = "maa://22524 (94.71), *maa://22432 (77.78)"
$ws.Range("X22").Value = "maa://21282 (98.51), *maa://37649 (67.86)"
$ws.Range("L23").Value = "maa://39756 (94.7), maa://39875 (93.85)"
$ws.Range("AB26").Value = "maa://42235 (93.9)"
$ws.Range("X28").Value = "maa://39929 (90.15), maa://41749 (90.77), ***maa://39723 (14.29)"
$ws.Range("P29").Value = "*maa://23168 (57.14), *maa://30050 (51.61)"
$ws.Range("AB30").Value = "maa://42979 (96.03), maa://45045 (100.0)"
$ws.Range("T32").Value = "maa://42859 (95.74), maa://41108 (88.0), maa://41238 (96.55)"
$ws.Range("AF32").Value = "*maa://42408 (80.0)"
$ws.Range("L35").Value = "maa://41296 (96.9)"
$ws.Range("T36").Value = "maa://27613 (99.05)"
$ws.Range("T37").Value = "**maa://39354 (33.33)"
$ws.Range("H39").Value = "maa://25199 (84.82), maa://36670 (87.91), maa://30434 (90.14), ***maa://25036 (16.0), *maa://44165 (66.67), *maa://45059 (66.67)"
$ws.Range("O40").NumberFormat = "@"
$ws.Range("O40").Value = "4"
$ws.Range("P40").Value = "maa://23278 (95.7), maa://21386 (95.74), maa://36664 (90.91), maa://45550 (100.0)"
$ws.Range("P41").Value = "**maa://35616 (38.24), maa://43177 (87.5)"
$ws.Range("T44").Value = "maa://39366 (87.5)"
$ws.Range("T45").Value = "**maa://39364 (38.46)"
$ws.Range("H60").Value = "*maa://40438 (64.71)"
